# Applies the "Boolean Do Suppliers Bid at Peak Capacity Factors" update:
#  - About sheet: replace the California-specific notes block with the
#    generic national "About" text (Certain plant types.../When bidding.../
#    This variable helps... paragraphs), dropping the old CA-only notes.
#  - BDSBaPCF sheet: rename the header, flip several booleans, and add
#    crude oil / heavy-or-residual-fuel-oil / municipal-solid-waste rows.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("About")
$ws2 = $wb.Worksheets.Item("BDSBaPCF")

# ---------------------------------------------------------------------
# 1) About sheet - rebuild the notes area from scratch.
# ---------------------------------------------------------------------
$ws1.Range("A1:B28").EntireRow.Delete()

$ws1.Range("A1").Value = "BDSBaPCF Boolean Do Suppliers Bid at Peak Capacity Factors"
$ws1.Range("A1").Font.Bold = $true

$ws1.Range("A3").Value = "Source:"
$ws1.Range("A3").Font.Bold = $true
$ws1.Range("B3").Value = "None needed.  See notes below."

$ws1.Range("A5").Value = "Notes"
$ws1.Range("A5").Font.Bold = $true

$ws1.Range("A6").Value = "Certain plant types, such as coal and natural gas, are capable of running for most"
$ws1.Range("A7").Value = "of the year (at their Peak Time capacity factors), if there were demand for their"
$ws1.Range("A8").Value = "services.  Other plant types, such as solar, wind, and hydro, are limited to a"
$ws1.Range("A9").Value = "much smaller fraction of the year by physical constraints (such as lack of sunlight"
$ws1.Range("A10").Value = "during some hours)."

$ws1.Range("A12").Value = "When bidding for how much power they are willing to supply, plant operators"
$ws1.Range("A13").Value = "that can supply more power are likely to bid it, since they will make more money"
$ws1.Range("A14").Value = "if they are able to run their coal or even a relatively inefficient peaker plant"
$ws1.Range("A15").Value = "more often.  (As long as the plant is built, they want to get use out of it.)"

$ws1.Range("A17").Value = "This variable helps the model decide which capacity factors best represent"
$ws1.Range("A18").Value = "the ones a plant owner would bid (offer to sell), which guides the dispatch"
$ws1.Range("A19").Value = "mechanism.  Plants that could conceivably be run full-out all year (minus"
$ws1.Range("A20").Value = "required maintenance, etc.) bid at their Peak Capacity Factors;"
$ws1.Range("A21").Value = "other plants bid at their Expected Capacity Factors, which account for all"
$ws1.Range("A22").Value = "of the things that prevent a plant from running at its peak all the time."

$ws1.Range("A1:A22").EntireRow.AutoFit()

# ---------------------------------------------------------------------
# 2) BDSBaPCF sheet - header text + wrap, boolean flips, new fuel rows.
# ---------------------------------------------------------------------
$ws2.Range("B1").Value = "Do Suppliers Bid at Peak Capacity Factors (Boolean)"
$ws2.Range("B1").Font.Bold = $true
$ws2.Range("B1").WrapText = $true
$ws2.Rows.Item(1).RowHeight = 45

$ws2.Range("B2").Value = 1
$ws2.Range("B9").Value = 1
$ws2.Range("B12").Value = 0

$ws2.Range("A15").Value = "crude oil"
$ws2.Range("B15").Formula = "=B11"

$ws2.Range("A16").Value = "heavy or residual fuel oil"
$ws2.Range("B16").Formula = "=B11"

$ws2.Range("A17").Value = "municipal solid waste"
$ws2.Range("B17").Formula = "=B9"

$ws2.Range("A2:A17").EntireRow.AutoFit()

# ---------------------------------------------------------------------
# 3) Active sheet / selection bookkeeping.
# ---------------------------------------------------------------------
[void]$ws2.Range("B5").Select()
[void]$ws1.Activate()
[void]$ws1.Range("A1").Select()
